$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, derived from the authoritative diff.
$updates = @{
    "D2" = "299.32"
    "E2" = "-1.94%"
    "E3" = "-1.56%"
    "D4" = "5.105"
    "E4" = "-1.48%"
    "D5" = "0.07955"
    "E5" = "5.77%"
    "D6" = "2.283"
    "E6" = "-1.86%"
    "D7" = "7.771"
    "E7" = "-3.12%"
    "D8" = "3.865"
    "E8" = "-0.38%"
    "D9" = "0.9262"
    "E9" = "1.02%"
    "D10" = "0.1736"
    "E10" = "-0.26%"
    "D11" = "0.07540"
    "E11" = "-0.75%"
    "D12" = "0.09359"
    "E12" = "13.15%"
    "D13" = "0.03046"
    "E13" = "0.09%"
    "D14" = "0.1004"
    "E14" = "1.09%"
    "D15" = "0.001510"
    "E15" = "0.04%"
    "D16" = "0.005920"
    "E16" = "-2.34%"
    "D17" = "3.483"
    "E17" = "-0.33%"
    "E18" = "1.46%"
    "E19" = "0.25%"
    "E20" = "0.44%"
    "D21" = "3.924"
    "E21" = "-15.61%"
    "D22" = "0.1701"
    "E22" = "8.73%"
    "D23" = "0.04610"
    "E23" = "-0.30%"
    "D24" = "0.001250"
    "E24" = "-0.94%"
    "D25" = "0.004482"
    "E25" = "-1.18%"
    "E26" = "-7.65%"
    "D27" = "0.0003398"
    "E27" = "24.02%"
    "D39" = "0.01756"
    "E39" = "-0.01%"
    "D40" = "0.04628"
    "E40" = "1.14%"
    "D41" = "0.006958"
    "E41" = "-3.74%"
    "D42" = "0.1362"
    "E42" = "-0.24%"
    "E43" = "-0.42%"
    "D44" = "0.01024"
    "E44" = "-5.01%"
    "D45" = "0.00006274"
    "E45" = "-3.81%"
    "D46" = "0.00000000750"
    "E46" = "-0.01%"
    "D47" = "0.007980"
    "E47" = "-19.31%"
    "D48" = "1.157"
    "E48" = "40.97%"
    "D49" = "0.00002100"
    "E49" = "-0.01%"
    "D50" = "0.0002000"
    "E50" = "0.06%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (prices, percentages)
    # are preserved exactly as literal text, matching the source data.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
